# Fruta / hortaliza, semanal
# A new weekly price record (Damasco, Castle Brite, Primera, Paine) is added
# at the top of the data (row 20), pushing the existing rows 20-46 down to
# rows 21-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 20; this shifts rows 20-46 down to 21-47
# and extends the used range to A1:T47 automatically.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(20, 1).Value = 6
$ws.Cells.Item(20, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44519
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100103
$ws.Cells.Item(20, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(20, 9).Value = 100103003
$ws.Cells.Item(20, 10).Value = "Damasco"
$ws.Cells.Item(20, 11).Value = "Castle Brite"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 150
$ws.Cells.Item(20, 14).Value = 22000
$ws.Cells.Item(20, 15).Value = 24000
$ws.Cells.Item(20, 16).Value = 23000
$ws.Cells.Item(20, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(20, 18).Value = "Paine"
$ws.Cells.Item(20, 19).Value = 1438
$ws.Cells.Item(20, 20).Value = 16
